$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.893.68'
$ws.Range('E2').Value = '  -0.46%  '

$ws.Range('D3').Value = '1.810.05'
$ws.Range('E3').Value = '  +1.18%  '

$ws.Range('D4').Value = '0.9960'
$ws.Range('E4').Value = '  -0.95%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '307.85'
$ws.Range('E5').Value = '  -1.85%  '

$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').Value = '0.9988'
$ws.Range('E6').Value = '  -0.53%  '

$ws.Range('D7').Value = '0.4997'
$ws.Range('E7').Value = '  -4.06%  '

$ws.Range('D8').Value = '0.3882'
$ws.Range('E8').Value = '  +2.26%  '

$ws.Range('D9').Value = '0.09406'
$ws.Range('E9').Value = '  +19.83%  '

$ws.Range('D10').Value = '1.098'
$ws.Range('E10').Value = '  +0.84%  '

$ws.Range('D11').Value = '40.52'
$ws.Range('E11').Value = '  -1.99%  '

$ws.Range('D12').Value = '6.310'
$ws.Range('E12').Value = '  +0.96%  '

$ws.Range('D13').Value = '1.000'
$ws.Range('E13').Value = '  -0.54%  '

$ws.Range('D14').Value = '20.56'
$ws.Range('E14').Value = '  +0.73%  '

$ws.Range('D15').Value = '1.793.34'
$ws.Range('E15').Value = '  -0.08%  '

$ws.Range('D16').Value = '7.213'
$ws.Range('E16').Value = '  -0.54%  '

$ws.Range('D17').Value = '0.00001123'
$ws.Range('E17').Value = '  +4.14%  '

$ws.Range('D18').Value = '92.73'
$ws.Range('E18').Value = '  +1.19%  '

$ws.Range('D19').Value = '0.06579'
$ws.Range('E19').Value = '  +0.74%  '

$ws.Range('D20').Value = '0.9994'
$ws.Range('E20').Value = '  -0.48%  '

$ws.Range('D21').Value = '17.11'
$ws.Range('E21').Value = '  -0.57%  '

$ws.Range('D22').Value = '5.924'
$ws.Range('E22').Value = '  -0.14%  '

$ws.Range('D23').Value = '27.848.67'
$ws.Range('E23').Value = '  -0.72%  '

$ws.Range('D24').Value = '11.00'
$ws.Range('E24').Value = '  -0.63%  '

$ws.Range('D25').Value = '2.227'
$ws.Range('E25').Value = '  -1.47%  '

$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '20.67'
$ws.Range('E26').Value = '  +1.73%  '

$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '156.40'
$ws.Range('E27').Value = '  -2.82%  '

$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '2.421'
$ws.Range('E28').Value = '  +4.90%  '

$ws.Range('B29').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C29').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D29').Value = '2.003.75'
$ws.Range('E29').Value = '  +0.50%  '

$ws.Range('D30').Value = '127.28'
$ws.Range('E30').Value = '  +4.13%  '

$ws.Range('D31').Value = '0.1071'
$ws.Range('E31').Value = '  +0.27%  '

$ws.Range('D32').Value = '1.058'
$ws.Range('E32').Value = '  +1.43%  '

$ws.Range('D33').Value = '5.550'
$ws.Range('E33').Value = '  +0.89%  '

$ws.Range('D34').Value = '3.606'
$ws.Range('E34').Value = '  -1.93%  '

$ws.Range('D35').Value = '0.06815'
$ws.Range('E35').Value = '  -5.26%  '

$ws.Range('D36').Value = '8.888'
$ws.Range('E36').Value = '  +1.73%  '

$ws.Range('D37').Value = '0.02301'
$ws.Range('E37').Value = '  -0.37%  '

$ws.Range('D38').Value = '0.2140'
$ws.Range('E38').Value = '  +0.73%  '

$ws.Range('D39').Value = '11.39'
$ws.Range('E39').Value = '  -5.98%  '

$ws.Range('D40').Value = '4.926'
$ws.Range('E40').Value = '  -2.18%  '

$ws.Range('D41').Value = '0.6169'
$ws.Range('E41').Value = '  +1.06%  '

$ws.Range('E42').Value = '  -0.29%  '

$ws.Range('D43').Value = '1.140'
$ws.Range('E43').Value = '  -1.79%  '

$ws.Range('D44').Value = '13.05'
$ws.Range('E44').Value = '  -1.23%  '

$ws.Range('D45').Value = '0.5875'
$ws.Range('E45').Value = '  -0.42%  '

$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '3.659'
$ws.Range('E46').Value = '  -2.59%  '

$ws.Range('B47').Value = 'WEMIXTOKEN'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '1.271'
$ws.Range('E47').Value = '  -6.80%  '

$ws.Range('D48').Value = '123.56'
$ws.Range('E48').Value = '  -3.25%  '

$ws.Range('D49').Value = '1.936'
$ws.Range('E49').Value = '  +1.54%  '

$ws.Range('D50').Value = '1.174'
$ws.Range('E50').Value = '  -3.87%  '

$ws.Range('D51').Value = '0.06713'
$ws.Range('E51').Value = '  -0.19%  '
